$wb = $excel.ActiveWorkbook

# The "二级科室" (Secondary Department) column is being removed from the
# "02" worksheet (sheet4.xml / rId4). Deleting the entire column D shifts
# all subsequent columns left by one, which matches the target diff
# (dimension A1:O4 -> A1:N4, column widths/styles shifted, shared string
# "二级科室" removed).
$ws = $wb.Worksheets.Item("02")
$ws.Range("D:D").Delete()
